$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New predicted values for columns B and C (rows 2-67),
# resulting from re-running the prediction loop with different outcomes
# than the previous run (per commit message).
$newValues = @(
    @(2, -0.2905807551338088, -0.2754734166150056),
    @(3, 0.2520100156392277, 0.2512509213060825),
    @(4, 0.181772945154432, 0.1904038857545432),
    @(5, -0.03464747234919305, -0.03003710391555485),
    @(6, 0.182808586937076, 0.1898816109089304),
    @(7, -0.4270243578396969, -0.4176839069414724),
    @(8, -0.276055745960919, -0.2517999141463262),
    @(9, -0.3691123619183683, -0.3682023173947372),
    @(10, 0.3475007893562577, 0.358590003820514),
    @(11, -0.1948147226932498, -0.1873195062811653),
    @(12, -0.04757160104148836, -0.0366425372575428),
    @(13, 0.007943755213636949, 0.006725963445144187),
    @(14, 0.1190593680963605, 0.1255110335824472),
    @(15, 0.005615524014902895, 0.01632118054256226),
    @(16, 0.382792554714191, 0.3979897967195563),
    @(17, 0.663893829496523, 0.6550865373073558),
    @(18, 0.08809615446765026, 0.05402683447005929),
    @(19, 0.4392482395985298, 0.430794098429519),
    @(20, 0.3442132110085145, 0.3548582582374352),
    @(21, 0.5114961857042081, 0.5454427521177849),
    @(22, 0.4268697929224574, 0.4265546718250002),
    @(23, -0.009929473126293969, -0.005536045819890062),
    @(24, 4.312011759827991, 4.299034319154053),
    @(25, 0.4963507696657999, 0.468139487400885),
    @(26, 0.4292720604632921, 0.4053728865731186),
    @(27, 0.2748891192003126, 0.2557785342226873),
    @(28, 1.029497876602861, 0.9990042108057645),
    @(29, 5.689326723106623, 5.194547144540082),
    @(30, 0.972615419403221, 0.9269082234680981),
    @(31, -0.2284627189200214, -0.2640449344501448),
    @(32, 0.781414073567569, 0.7495914205171905),
    @(33, 0.8716976586112575, 0.855107779250092),
    @(34, -0.6712122212930119, -0.6845279241454052),
    @(35, 0.8104279385953042, 0.8094908584298817),
    @(36, 0.7520278950017755, 0.7506866628945279),
    @(37, 0.732739764398458, 0.7275018373043871),
    @(38, 0.7430073150876311, 0.7260916723852486),
    @(39, 0.5631300873331944, 0.5697552768097492),
    @(40, 0.7318677309714507, 0.7401914074990478),
    @(41, 0.5660711582779037, 0.56256867322948),
    @(42, 0.6848642555562446, 0.671167544775667),
    @(43, 0.7177917397137505, 0.7088763397019342),
    @(44, 0.6558159246968844, 0.6577395631667096),
    @(45, 0.6138107185281998, 0.6187164044392784),
    @(46, -1.276305821483627, -1.268471728031601),
    @(47, -0.9906637592918169, -0.9850208239613125),
    @(48, -0.8865357028131678, -0.8794563246545962),
    @(49, -0.649102846523711, -0.6423910600974911),
    @(50, -0.05857146424311572, -0.05412765171008823),
    @(51, -0.8755976832279182, -0.8661956392887654),
    @(52, -0.8755976832279182, -0.8661956392887654),
    @(53, -1.106815814517363, -1.106885937752893),
    @(54, -0.1983086422808562, -0.1898439127467212),
    @(55, -1.012619185558019, -1.004814373567967),
    @(56, -0.910328116220338, -0.8957979785958557),
    @(57, -0.9555973097622447, -0.94207083815229),
    @(58, -1.156262168401689, -1.129295974554265),
    @(59, -0.8753434166570127, -0.8527335191660237),
    @(60, -0.5126573540511948, -0.4890907932833319),
    @(61, 0.3666079006070712, 0.3705108953312143),
    @(62, -1.246597742964313, -1.226015649819886),
    @(63, -0.7596008310974316, -0.7209243555231546),
    @(64, -0.8947759556778864, -0.8856412595839435),
    @(65, -0.1334905590660783, -0.1091942891879161),
    @(66, -0.8243742046194242, -0.7979388346493719),
    @(67, -0.8273268430114561, -0.7904733825880573),
)

foreach ($entry in $newValues) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
